$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 101, shifting existing rows 101:161 down to 102:162
$ws.Rows("101:101").Insert()

# Populate the newly inserted row 101 with the new record's data.
# Columns A,B,C,E,F,G,H,I,O,R keep the same values as the rest of the
# "Berenjena" / "Terminal La Palmera de La Serena" dataset rows.
$ws.Range("A101").Value = 8
$ws.Range("B101").Value = "Terminal La Palmera de La Serena"
$ws.Range("C101").Value = "Coquimbo"
$ws.Range("D101").Value = 44813
$ws.Range("E101").Value = 4
$ws.Range("F101").Value = 100112001
$ws.Range("G101").Value = "Berenjena"
$ws.Range("H101").Value = "Sin especificar"
$ws.Range("I101").Value = "Primera"
$ws.Range("J101").Value = 560
$ws.Range("K101").Value = 10000
$ws.Range("L101").Value = 11000
$ws.Range("M101").Value = 10500
$ws.Range("N101").Value = "$/caja 40 unidades"
$ws.Range("O101").Value = "Región de Arica y Parinacota"
$ws.Range("P101").Value = 262
$ws.Range("Q101").Value = 40
$ws.Range("R101").Value = "Hortaliza"
